# Add files via upload
# Additional Codes Of "Buy Ticket" Section
#
# The header row of the fixture sheet is retyped with proper
# capitalisation / spacing, and the active selection is moved to F1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = " Team 1"
$ws.Range("C1").Value = " Team 2"
$ws.Range("D1").Value = " Time"
$ws.Range("E1").Value = "Venue"
$ws.Range("F1").Value = "Check"

$ws.Range("F1").Select()
